# Insert a new data row at row 8 (pushing existing rows 8..126 down to 9..127)
# and populate the new row 8 with the new daily price-record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8..126 down by one row, creating a new, empty row 8.
$ws.Rows.Item(8).Insert()

# Fill the new row 8 with the new record (same categorical data as the
# surrounding rows, new date/volume, same min/max/avg price and $/kg).
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 45169
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100108
$ws.Range("H8").Value = "Tropicales y subtropicales"
$ws.Range("I8").Value = 100108007
$ws.Range("J8").Value = "Coco"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 25
$ws.Range("N8").Value = 36000
$ws.Range("O8").Value = 36000
$ws.Range("P8").Value = 36000
$ws.Range("Q8").Value = "$/malla 20 unidades"
$ws.Range("R8").Value = "Perú"
$ws.Range("S8").Value = 1800
$ws.Range("T8").Value = 20
